# Corrige resultados experimentos AQ
# Updates the AQ_results input data (columns B-G) to corrected values.
# Formula cells (I-N, and summary rows 38-44) recalc automatically.
# Also restores the intended view state: AQ_results active/selected,
# both sheets zoomed to 85%, and their last-used selections.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("AQ_results")
$ws2 = $wb.Worksheets.Item("AQ_descrip")

$ws1.Range("D2").Value = -0.0446439316594496
$ws1.Range("B3").Value = 0
$ws1.Range("C3").Value = -0.0700835299053749
$ws1.Range("D3").Value = 0
$ws1.Range("E3").Value = -0.0700835299053749
$ws1.Range("F3").Value = -0.148531108575075
$ws1.Range("G3").Value = -0.0239034341045996
$ws1.Range("D4").Value = -0.0087643617075373
$ws1.Range("D5").Value = 0.108013138900988
$ws1.Range("D6").Value = 0.193677566673
$ws1.Range("D7").Value = -0.150555304547262
$ws1.Range("D8").Value = -0.163737569060625
$ws1.Range("D9").Value = 0.0124432309833
$ws1.Range("D11").Value = -0.0454686432291996
$ws1.Range("D12").Value = 0
$ws1.Range("B13").Value = -0.039822718860175
$ws1.Range("C13").Value = 0.0115011149338753
$ws1.Range("D13").Value = -0.039822718860175
$ws1.Range("E13").Value = 0
$ws1.Range("F13").Value = -0.0433134937550496
$ws1.Range("G13").Value = -0.0538636614733498
$ws1.Range("D14").Value = -0.00188475789767455
$ws1.Range("D15").Value = 0.0134233426732749
$ws1.Range("D16").Value = -0.0095806796875246
$ws1.Range("D17").Value = 0.00792748446512503
$ws1.Range("D18").Value = 0.022819650385625
$ws1.Range("D19").Value = 0.00618741288557496
$ws1.Range("B20").Value = -0.0578137289112999
$ws1.Range("C20").Value = 0.02916874961435
$ws1.Range("D20").Value = -0.0578137289112999
$ws1.Range("E20").Value = 0.02916874961435
$ws1.Range("F20").Value = -0.0369035856321249
$ws1.Range("G20").Value = 0.0410770820361752
$ws1.Range("D21").Value = -0.0354782550030247
$ws1.Range("B22").Value = 0.0271867755909252
$ws1.Range("C22").Value = 0.0181186998998251
$ws1.Range("D22").Value = 0.0354936559695254
$ws1.Range("E22").Value = 0.0260380378004752
$ws1.Range("F22").Value = 0.0148573387720252
$ws1.Range("G22").Value = 0.0157138250000002
$ws1.Range("D23").Value = -0.00357938964344956
$ws1.Range("B24").Value = 0
$ws1.Range("C24").Value = 0
$ws1.Range("D24").Value = 0
$ws1.Range("E24").Value = 0
$ws1.Range("F24").Value = 0
$ws1.Range("G24").Value = 0
$ws1.Range("B25").Value = 0.0342201545823253
$ws1.Range("C25").Value = 0.0223623298828002
$ws1.Range("D25").Value = -0.0624635318309745
$ws1.Range("E25").Value = -0.0735452979467498
$ws1.Range("F25").Value = 0.0217203643579751
$ws1.Range("G25").Value = 0.0132804041015502
$ws1.Range("D26").Value = -0.0265632264823748
$ws1.Range("D27").Value = -0.00433576854472503
$ws1.Range("D29").Value = -0.01291382490985
$ws1.Range("D30").Value = 0.00157834892827511
$ws1.Range("B31").Value = 0
$ws1.Range("C31").Value = 0.015399626657675
$ws1.Range("D31").Value = 0
$ws1.Range("E31").Value = 0.0153885834886001
$ws1.Range("F31").Value = 0.0153885834886001
$ws1.Range("G31").Value = 0
$ws1.Range("D32").Value = 0.00277195898437518
$ws1.Range("D33").Value = 0.00842755801782516
$ws1.Range("D34").Value = 0.00155555341047501
$ws1.Range("B35").Value = -0.0913916513671748
$ws1.Range("C35").Value = -0.0951902273487749
$ws1.Range("D35").Value = -0.0923554627904748
$ws1.Range("E35").Value = -0.0951902273487749
$ws1.Range("F35").Value = -0.0788742201122247
$ws1.Range("G35").Value = -0.0988360378856248
$ws1.Range("D36").Value = -0.0299966226562747


# View state: set AQ_descrip's selection/zoom first, then finish on
# AQ_results so it ends up the active/selected sheet (activeTab=0).
$ws2.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 85
$ws2.Range("A1").Select() | Out-Null

$ws1.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 85
$ws1.Range("A36").Select() | Out-Null
